$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of J column
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Row 14-17: summary stats with labels in column A and values in column B
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style B14 - bold, size 12, vertical center - then propagate via copy/paste
# formats so we don't mint extra intermediate cell styles for the range.
$b14 = $ws.Range("B14")
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.VerticalAlignment = -4108

$b14.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row heights for the new summary rows
$ws.Range("A14:A17").RowHeight = 15.6

# Page setup (paper size / orientation) as in the diff
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection as in the diff
$ws.Range("A14:B17").Select() | Out-Null
